$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 (Subj header row) values for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) meanEMG / legmaxROM values for columns B:E
$ws.Range("B2").Value = 60.780713310805943
$ws.Range("C2").Value = 49.183925646657009
$ws.Range("D2").Value = 64.887966098158202
$ws.Range("E2").Value = 52.598329430698286

# Update row 3 (STR) meanEMG / legmaxROM values for columns B:E
$ws.Range("B3").Value = 63.221206623705854
$ws.Range("C3").Value = 44.659973050356776
$ws.Range("D3").Value = 70.294355396639375
$ws.Range("E3").Value = 44.245052626267544

# Match the updated selection left behind in the saved worksheet view
$ws.Range("B1:E3").Select() | Out-Null
